$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns before writing, so that
# numeric-looking strings (e.g. "583.20", "1.00") are preserved verbatim
# as text instead of being auto-converted to numbers by Excel.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "60.504.88"
$ws.Range("E2").Value = "  +2.30%  "

$ws.Range("D3").Value = "2.600.90"
$ws.Range("E3").Value = "  +1.51%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "583.20"
$ws.Range("E5").Value = "  +6.13%  "

$ws.Range("D6").Value = "143.45"
$ws.Range("E6").Value = "  +2.53%  "

$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +1.26%  "

$ws.Range("D9").Value = "2.628.01"
$ws.Range("E9").Value = "  +2.51%  "

$ws.Range("D10").Value = "6.55"
$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("E11").Value = "  +2.52%  "

$ws.Range("E12").Value = "  -3.30%  "

$ws.Range("D13").Value = "0.372"
$ws.Range("E13").Value = "  +5.84%  "

$ws.Range("D14").Value = "3.086.66"
$ws.Range("E14").Value = "  +2.17%  "

$ws.Range("D15").Value = "24.54"
$ws.Range("E15").Value = "  +6.76%  "

$ws.Range("D16").Value = "60.538.95"
$ws.Range("E16").Value = "  +2.32%  "

$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").Value = "  +3.93%  "

$ws.Range("D18").Value = "2.625.71"
$ws.Range("E18").Value = "  +1.91%  "

$ws.Range("D19").Value = "11.38"
$ws.Range("E19").Value = "  +11.12%  "

$ws.Range("D20").Value = "4.68"
$ws.Range("E20").Value = "  +3.55%  "

$ws.Range("D21").Value = "348.28"
$ws.Range("E21").Value = "  +3.60%  "

$ws.Range("D22").Value = "6.91"
$ws.Range("E22").Value = "  +8.18%  "

$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("E24").Value = "  +9.55%  "

$ws.Range("D25").Value = "63.10"
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("D28").Value = "7.94"
$ws.Range("E28").Value = "  +7.85%  "

$ws.Range("D29").Value = "0.0₃0798"
$ws.Range("E29").Value = "  +4.83%  "

$ws.Range("D30").Value = "1.87"
$ws.Range("E30").Value = "  +12.26%  "

$ws.Range("D31").Value = "6.38"
$ws.Range("E31").Value = "  +4.10%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").Value = "164.51"
$ws.Range("E33").Value = "  +3.82%  "

$ws.Range("E34").Value = "  +2.83%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  +4.62%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +12.77%  "

$ws.Range("E37").Value = "  +6.58%  "

$ws.Range("E38").Value = "  +10.93%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "37.96"
$ws.Range("E39").Value = "  +1.51%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "313.43"
$ws.Range("E40").Value = "  +9.91%  "

$ws.Range("D41").Value = "3.89"
$ws.Range("E41").Value = "  +6.72%  "

$ws.Range("D42").Value = "0.842"
$ws.Range("E42").Value = "  -0.30%  "

$ws.Range("D43").Value = "135.20"
$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "5.09"
$ws.Range("E44").Value = "  +13.81%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.0991"
$ws.Range("E45").Value = "  +2.46%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "0.996"
$ws.Range("E46").Value = "  -0.31%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "19.87"
$ws.Range("E47").Value = "  +6.49%  "

$ws.Range("E48").Value = "  +3.38%  "

$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0553"
$ws.Range("E49").Value = "  +4.93%  "

$ws.Range("D50").Value = "20.23"
$ws.Range("E50").Value = "  +9.42%  "

$ws.Range("D51").Value = "0.0242"
$ws.Range("E51").Value = "  +4.41%  "

# Restore the default cell style so no stray number-format styling is
# left attached to the cells (matches original workbook formatting).
$priceVolRange.Style = "Normal"
